# Updates crypto price/volume data in the worksheet to reflect the latest
# GitHub Actions scrape (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.308.63"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.150.56"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'228.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").Value = "'62.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'15.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.00%  "
$ws.Range("D13").Value = "2.472.20"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").Value = "'22.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "'5.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "2.142.79"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "39.432.29"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'72.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'6.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "'228.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").Value = "'171.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'19.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("D30").Value = "'1.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  +10.15%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'4.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'7.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.71%  "
$ws.Range("D36").Value = "'0.0620"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'18.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").Value = "'102.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("D43").Value = "1.535.06"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  +6.52%  "
$ws.Range("E45").Value = "  +7.19%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'4.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "'2.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").Value = "2.356.79"
$ws.Range("E51").Value = "  +3.23%  "

Write-Host "Applied 94 cell updates"
